# Update status text "Ready for handoff" -> "In Translation" across all sheets,
# and narrow the "Status" column width on each sheet (columns E/F on Overview,
# column C on zh-cn and de-de) from 17.2159881591797 to 13.4101845877511.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status values.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the previously-widened status columns back down. Excel's ColumnWidth
# property is expressed in "characters" and gets snapped to the workbook's
# default-font pixel grid on save, so 12.5 is the input that lands closest to
# the target stored width (~13.41 character-units).
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
